$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 8.450800000000001
$ws.Range("A9").Value = -22.35589999999999
$ws.Range("B9").Value = 4.501400000000005
$ws.Range("D9").Value = -7.986099999999998
$ws.Range("B11").Value = 5.151300000000001
$ws.Range("A13").Value = -22.32200000000001
$ws.Range("A16").Value = -21.55629999999999
$ws.Range("B16").Value = 4.886299999999998
$ws.Range("A18").Value = -22.1112
$ws.Range("A20").Value = -21.65009999999999
$ws.Range("D22").Value = -7.945999999999999
$ws.Range("B23").Value = 7.992300000000008
$ws.Range("B24").Value = 4.479900000000003
$ws.Range("A26").Value = -21.06089999999997
$ws.Range("B26").Value = 5.976700000000002
$ws.Range("A27").Value = -22.00099999999999
$ws.Range("D27").Value = -8.280500000000004
$ws.Range("A29").Value = -21.04599999999999
$ws.Range("D29").Value = -7.886899999999996
$ws.Range("D32").Value = -7.199099999999991
$ws.Range("B34").Value = 9.4793
$ws.Range("A35").Value = -21.8039
$ws.Range("B35").Value = 4.7526
$ws.Range("A36").Value = -21.6267
$ws.Range("D37").Value = -7.117900000000001
$ws.Range("D38").Value = -7.412299999999997
$ws.Range("D39").Value = -7.669399999999996
$ws.Range("D41").Value = -7.855100000000004
$ws.Range("B44").Value = 4.609100000000005
$ws.Range("A45").Value = -21.60059999999998
$ws.Range("D45").Value = -7.255499999999997
$ws.Range("B48").Value = 4.909400000000006
$ws.Range("D48").Value = -8.205100000000002
$ws.Range("B49").Value = 5.6339
$ws.Range("D51").Value = -8.675499999999996
$ws.Range("B52").Value = 5.384499999999997
$ws.Range("A55").Value = -22.19700000000001
$ws.Range("D56").Value = -8.881500000000003
$ws.Range("A57").Value = -21.96760000000001
$ws.Range("D57").Value = -8.471799999999996
$ws.Range("D61").Value = -8.091800000000001
$ws.Range("D64").Value = -7.18589999999999
$ws.Range("B66").Value = 4.860299999999996
$ws.Range("B67").Value = 5.502800000000002
$ws.Range("A69").Value = -21.5597
$ws.Range("B73").Value = 9.292899999999996
$ws.Range("D75").Value = -8.225499999999998
$ws.Range("A76").Value = -19.52459999999997
$ws.Range("A78").Value = -21.72960000000001
$ws.Range("B78").Value = 5.557500000000001
$ws.Range("B80").Value = 9.523899999999994
$ws.Range("A82").Value = -22.06140000000002
$ws.Range("D82").Value = -8.789099999999996
$ws.Range("A83").Value = -21.5448
$ws.Range("D90").Value = -7.432299999999993
$ws.Range("B91").Value = 4.9665
$ws.Range("A93").Value = -21.41030000000001
$ws.Range("D93").Value = -6.91249999999999
$ws.Range("A97").Value = -21.5735
$ws.Range("B97").Value = 5.0064
$ws.Range("B99").Value = 5.6863
$ws.Range("D102").Value = -7.512800000000002
$ws.Range("B104").Value = 9.975200000000003
$ws.Range("D105").Value = -7.752500000000002
